$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 82.83048866666665
$ws.Range("H2").Value = 248.491466
$ws.Range("I2").Value = 0.3167437020391103
$ws.Range("J2").Value = 0.3167437020391103
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 409.6166503333334
$ws.Range("N2").Value = 1228.849951
$ws.Range("O2").Value = 0.6234125531262766
$ws.Range("P2").Value = 0.6234125531262766
$ws.Range("Q2").Value = 33928.74731311313
$ws.Range("R2").Value = 305358.7258180182
$ws.Range("S2").Value = 0.1974619999748704
$ws.Range("T2").Value = 0.1974619999748704

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 82.83048866666665
$ws.Range("H3").Value = 248.491466
$ws.Range("I3").Value = 0.3167437020391103
$ws.Range("J3").Value = 0.3167437020391103
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 56.495384
$ws.Range("N3").Value = 169.486152
$ws.Range("O3").Value = 0.08598266586728959
$ws.Range("P3").Value = 0.08598266586728959
$ws.Range("Q3").Value = 4679.540264130981
$ws.Range("R3").Value = 42115.86237717883
$ws.Range("S3").Value = 0.02723446789799715
$ws.Range("T3").Value = 0.02723446789799716

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 82.83048866666665
$ws.Range("H4").Value = 248.491466
$ws.Range("I4").Value = 0.3167437020391103
$ws.Range("J4").Value = 0.3167437020391103
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 190.9434713333333
$ws.Range("N4").Value = 572.830414
$ws.Range("O4").Value = 0.2906047810064339
$ws.Range("P4").Value = 0.2906047810064338
$ws.Range("Q4").Value = 15815.94103824966
$ws.Range("R4").Value = 142343.4693442469
$ws.Range("S4").Value = 0.09204723416624279
$ws.Range("T4").Value = 0.09204723416624279

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 152.851481
$ws.Range("H5").Value = 458.554443
$ws.Range("I5").Value = 0.5845039034954311
$ws.Range("J5").Value = 0.5845039034954312
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 409.6166503333334
$ws.Range("N5").Value = 1228.849951
$ws.Range("O5").Value = 0.6234125531262766
$ws.Range("P5").Value = 0.6234125531262766
$ws.Range("Q5").Value = 62610.51164570916
$ws.Range("R5").Value = 563494.6048113824
$ws.Range("S5").Value = 0.3643870707903615
$ws.Range("T5").Value = 0.3643870707903616

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 152.851481
$ws.Range("H6").Value = 458.554443
$ws.Range("I6").Value = 0.5845039034954311
$ws.Range("J6").Value = 0.5845039034954312
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 56.495384
$ws.Range("N6").Value = 169.486152
$ws.Range("O6").Value = 0.08598266586728959
$ws.Range("P6").Value = 0.08598266586728959
$ws.Range("Q6").Value = 8635.403114063705
$ws.Range("R6").Value = 77718.62802657334
$ws.Range("S6").Value = 0.05025720383237414
$ws.Range("T6").Value = 0.05025720383237414

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 152.851481
$ws.Range("H7").Value = 458.554443
$ws.Range("I7").Value = 0.5845039034954311
$ws.Range("J7").Value = 0.5845039034954312
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 190.9434713333333
$ws.Range("N7").Value = 572.830414
$ws.Range("O7").Value = 0.2906047810064339
$ws.Range("P7").Value = 0.2906047810064338
$ws.Range("Q7").Value = 29185.99238058105
$ws.Range("R7").Value = 262673.9314252294
$ws.Range("S7").Value = 0.1698596288726955
$ws.Range("T7").Value = 0.1698596288726955

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 25.824378
$ws.Range("H8").Value = 77.47313399999999
$ws.Range("I8").Value = 0.09875239446545848
$ws.Range("J8").Value = 0.0987523944654585
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 409.6166503333334
$ws.Range("N8").Value = 1228.849951
$ws.Range("O8").Value = 0.6234125531262766
$ws.Range("P8").Value = 0.6234125531262766
$ws.Range("Q8").Value = 10578.09521330183
$ws.Range("R8").Value = 95202.85691971643
$ws.Range("S8").Value = 0.06156348236104466
$ws.Range("T8").Value = 0.06156348236104468

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 25.824378
$ws.Range("H9").Value = 77.47313399999999
$ws.Range("I9").Value = 0.09875239446545848
$ws.Range("J9").Value = 0.0987523944654585
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 56.495384
$ws.Range("N9").Value = 169.486152
$ws.Range("O9").Value = 0.08598266586728959
$ws.Range("P9").Value = 0.08598266586728959
$ws.Range("Q9").Value = 1458.958151671152
$ws.Range("R9").Value = 13130.62336504037
$ws.Range("S9").Value = 0.008490994136918294
$ws.Range("T9").Value = 0.008490994136918296

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 25.824378
$ws.Range("H10").Value = 77.47313399999999
$ws.Range("I10").Value = 0.09875239446545848
$ws.Range("J10").Value = 0.0987523944654585
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 190.9434713333333
$ws.Range("N10").Value = 572.830414
$ws.Range("O10").Value = 0.2906047810064339
$ws.Range("P10").Value = 0.2906047810064338
$ws.Range("Q10").Value = 4930.996380344163
$ws.Range("R10").Value = 44378.96742309747
$ws.Range("S10").Value = 0.02869791796749554
$ws.Range("T10").Value = 0.02869791796749553

